$d = $word.ActiveDocument

# Replace "Introductory Application Development Concepts" with
# "Mobile Application Development" (collapsing the two runs into one).
$d.Content.Find.Execute("Introductory Application Development Concepts", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mobile Application Development", 2)
